$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 67
$ws.Range("I2").Value = 168
$ws.Range("J2").Value = 708
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 217
$ws.Range("N2").Value = 126
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 8
$ws.Range("S2").Value = 91
$ws.Range("T2").Value = 140
$ws.Range("V2").Value = 1123
$ws.Range("X2").Value = 1208
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 25
$ws.Range("AA2").Value = 14
